$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Furniture: purchases during the period reset to 0, cascading totals updated
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 3111072.93
$ws.Range("G4").Value = 1744051.14
$ws.Range("J4").Value = 174405.11
$ws.Range("L4").Value = 983382.76
$ws.Range("O4").Value = 758080.79

# Row 5 - Telecommunications: purchases during the period reset to 0, cascading totals updated
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 334523
$ws.Range("G5").Value = 334523
$ws.Range("J5").Value = 47789
$ws.Range("L5").Value = 284141.14
$ws.Range("O5").Value = 50381.83
